{"js": "// The template used Word \"field code\" syntax ( { m:userdoc 'zone1' } and\n// { m:enduserdoc } ) stored as real Word fields (fldChar begin/instrText/\n// fldChar end). The fix rewrites them as literal visible text runs\n// \"{\", \"m\", \":userdoc 'zone1'\", \"}\" (and \"{m:\" / \"enduserdoc}\") so the\n// M2Doc parser can tokenize them itself (TokenIteratorFieldRewriterSplit).\n//\n// We locate each Word field, capture its owning paragraph, delete the\n// field (removing the begin/instrText/end runs), then re-insert the\n// equivalent OOXML runs of plain text (<w:t>) in the same paragraph,\n// splitting them exactly like the target runs. For the second field the\n// paragraph also carries a `_GoBack` bookmark in between two runs, which\n// we recreate at the same spot.\n\nconst body = context.document.body;\n\nconst fields = body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length < 2) {\n  throw new Error(\"Expected 2 user-doc fields (zone1 / enduserdoc), found \" + fields.items.length);\n}\n\n// Use the paragraph collection to find the paragraphs that own each field,\n// by correlating field index before mutating anything.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find which paragraph currently contains each field's field-code text by\n// checking paragraph field counts (each target paragraph has exactly 1\n// field and nothing else).\nconst paraFieldCounts = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const pf = paragraphs.items[i].fields;\n  pf.load(\"items\");\n  paraFieldCounts.push(pf);\n}\nawait context.sync();\n\nlet zone1ParaIndex = -1;\nlet enduserdocParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const items = paraFieldCounts[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].load(\"code\");\n  }\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const items = paraFieldCounts[i].items;\n  for (let j = 0; j < items.length; j++) {\n    const code = items[j].code;\n    if (code.indexOf(\"userdoc\") !== -1 && code.indexOf(\"enduserdoc\") === -1) {\n      zone1ParaIndex = i;\n    } else if (code.indexOf(\"enduserdoc\") !== -1) {\n      enduserdocParaIndex = i;\n    }\n  }\n}\n\nif (zone1ParaIndex === -1 || enduserdocParaIndex === -1) {\n  throw new Error(\"Could not locate zone1/enduserdoc field paragraphs\");\n}\n\nconst zone1Para = paragraphs.items[zone1ParaIndex];\nconst enduserdocPara = paragraphs.items[enduserdocParaIndex];\n\nconst OOXML_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction pkg(innerParagraphXml) {\n  return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document ' + OOXML_NS + '><w:body>' + innerParagraphXml + '</w:body></w:document></pkg:xmlData>' +\n    '</pkg:part></pkg:package>';\n}\n\n// Delete the fields first (removes fldChar begin/instrText/fldChar end runs).\nzone1Para.fields.load(\"items\");\nenduserdocPara.fields.load(\"items\");\nawait context.sync();\nzone1Para.fields.items[0].delete();\nenduserdocPara.fields.items[0].delete();\nawait context.sync();\n\n// Re-insert plain-text runs, replacing the (now emptied) paragraph content.\nconst zone1Xml = pkg(\n  '<w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:userdoc \\'zone1\\'</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>'\n);\nzone1Para.getRange().insertOoxml(zone1Xml, Word.InsertLocation.replace);\n\nconst enduserdocXml = pkg(\n  '<w:p>' +\n    '<w:r><w:t>{m:</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">enduserdoc}</w:t></w:r>' +\n  '</w:p>'\n);\nenduserdocPara.getRange().insertOoxml(enduserdocXml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The template stores two M2Doc tags ( { m:userdoc 'zone1' } and\n# { m:enduserdoc } ) as real Word fields (fldChar begin/instrText/fldChar\n# end). The fix rewrites each one as literal visible text runs instead of\n# a field, so the M2Doc parser tokenizes the braces itself\n# (TokenIteratorFieldRewriterSplit):\n#   { m:userdoc 'zone1' }  ->  \"{\" \"m\" \":userdoc 'zone1'\" \"}\"\n#   { m:enduserdoc }       ->  \"{m:\" <bookmark> \"enduserdoc}\"\n#\n# Strategy: find each field by its instruction text, note which paragraph\n# (by index) owns it, delete the field (this removes the begin/instrText/\n# end runs but leaves the owning paragraph in place, now empty), then\n# re-insert the equivalent plain <w:t> runs via Range.InsertXML so the\n# run-splitting matches exactly. The second field's paragraph also has a\n# `_GoBack` bookmark sitting between two runs, which we recreate at the\n# same spot.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexForPosition($doc, $pos) {\n    $idx = 0\n    foreach ($p in $doc.Paragraphs) {\n        $idx++\n        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {\n            return $idx\n        }\n    }\n    return -1\n}\n\n$zone1Field = $null\n$enduserdocField = $null\nforeach ($f in $d.Fields) {\n    $code = $f.Code.Text\n    if ($code -like \"*enduserdoc*\") {\n        $enduserdocField = $f\n    } elseif ($code -like \"*userdoc*\") {\n        $zone1Field = $f\n    }\n}\n\nif ($null -eq $zone1Field -or $null -eq $enduserdocField) {\n    throw \"Could not locate zone1/enduserdoc fields\"\n}\n\n$zone1ParaIndex = Get-ParagraphIndexForPosition $d $zone1Field.Code.Start\n$enduserdocParaIndex = Get-ParagraphIndexForPosition $d $enduserdocField.Code.Start\n\nif ($zone1ParaIndex -eq -1 -or $enduserdocParaIndex -eq -1) {\n    throw \"Could not locate owning paragraphs for the fields\"\n}\n\n$zone1Field.Delete()\n$enduserdocField.Delete()\n\n$zone1XmlRun = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc ''zone1''</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$zone1Para = $d.Paragraphs.Item($zone1ParaIndex)\n$zone1Para.Range.InsertXML($zone1XmlRun)\n\n$enduserdocXmlRun = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{m:</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t xml:space=\"preserve\">enduserdoc}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$enduserdocPara = $d.Paragraphs.Item($enduserdocParaIndex)\n$enduserdocPara.Range.InsertXML($enduserdocXmlRun)\n\n$d.Save()\n"}
